# Apply the dated-worksheet update: refresh the date heading and all
# "two-digit ÷ one-digit" practice problems in the single table.

$d = $word.ActiveDocument

# --- Heading date -----------------------------------------------------
$d.Content.Find.Execute("2024-06-08 Saturday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-06-09 Sunday", 2) | Out-Null

# --- Table cell values --------------------------------------------------
# Mapping of (row, col) -> new text, addressed via the Table/Cell object
# model so cells are targeted unambiguously even where an old/new value
# collides with another cell's old/new value elsewhere in the table.
$tbl = $d.Tables.Item(1)

$cellUpdates = @(
    @{Row=1;  Col=1; New="46÷7="},
    @{Row=1;  Col=2; New="25÷9="},
    @{Row=1;  Col=3; New="30÷8="},
    @{Row=1;  Col=4; New="43÷4="},
    @{Row=1;  Col=5; New="44÷2="},

    @{Row=5;  Col=1; New="89÷8="},
    @{Row=5;  Col=2; New="15÷8="},
    @{Row=5;  Col=3; New="83÷9="},
    @{Row=5;  Col=4; New="93÷9="},
    @{Row=5;  Col=5; New="70÷5="},

    @{Row=9;  Col=1; New="50÷7="},
    @{Row=9;  Col=2; New="98÷5="},
    @{Row=9;  Col=3; New="46÷3="},
    @{Row=9;  Col=4; New="59÷7="},
    @{Row=9;  Col=5; New="25÷4="},

    @{Row=13; Col=1; New="72÷9="},
    @{Row=13; Col=2; New="59÷7="},
    @{Row=13; Col=3; New="13÷7="},
    @{Row=13; Col=4; New="42÷3="},
    @{Row=13; Col=5; New="31÷9="},

    @{Row=17; Col=1; New="29÷4="},
    @{Row=17; Col=2; New="40÷6="},
    @{Row=17; Col=3; New="30÷4="},
    @{Row=17; Col=4; New="71÷5="},
    @{Row=17; Col=5; New="42÷5="}
)

foreach ($u in $cellUpdates) {
    $cell = $tbl.Cell($u.Row, $u.Col)
    $rng = $cell.Range
    # Trim the trailing end-of-cell marker so only the visible text is
    # replaced, leaving paragraph/run formatting untouched.
    $rng.End = $rng.End - 1
    $rng.Text = $u.New
}
